# Auto-generated Excel COM-interop edit script
# Refreshes the crypto price/volume snapshot: updates Price (D) and Volume(1h) (E)
# values for the listed coins, and swaps the Kaspa / dogwifhat rows (38 & 39),
# matching the data pulled by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.929.41'
$ws.Range("E2").Value = '  -3.12%  '
$ws.Range("D3").Value = '2.919.39'
$ws.Range("E3").Value = '  -3.87%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.14'
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.65'
$ws.Range("E6").Value = '  -5.37%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -2.42%  '
$ws.Range("D9").Value = '2.918.23'
$ws.Range("E9").Value = '  -3.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.97'
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.145'
$ws.Range("E11").Value = '  -4.72%  '
$ws.Range("E12").Value = '  -3.89%  '
$ws.Range("E13").Value = '  -3.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.62'
$ws.Range("E14").Value = '  -5.76%  '
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").Value = '3.402.29'
$ws.Range("E16").Value = '  -3.88%  '
$ws.Range("D17").Value = '60.896.15'
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("E18").Value = '  -4.39%  '
$ws.Range("D19").Value = '2.920.28'
$ws.Range("E19").Value = '  -3.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '430.45'
$ws.Range("E20").Value = '  -5.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.63'
$ws.Range("E21").Value = '  -4.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.682'
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.14'
$ws.Range("E23").Value = '  -4.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.43'
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.82'
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  -3.35%  '
$ws.Range("E27").Value = '  -2.88%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -3.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.62'
$ws.Range("E31").Value = '  -3.03%  '
$ws.Range("E32").Value = '  -4.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.58'
$ws.Range("E33").Value = '  -3.81%  '
$ws.Range("E34").Value = '  -3.53%  '
$ws.Range("D35").Value = '0.0₃0872'
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("E36").Value = '  -2.69%  '
$ws.Range("E37").Value = '  -4.94%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.03'
$ws.Range("E38").Value = '  -5.80%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.128'
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.62'
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("E41").Value = '  -4.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.66'
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.297'
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.00'
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '377.77'
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("E46").Value = '  -3.15%  '
$ws.Range("D47").Value = '2.695.54'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.41'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.94'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("E51").Value = '  -1.89%  '
